# Daily attendance processing - 2025-12-31 13:46:09
# Normalises the "Recorded By" (column G) cell text: the comma-separated
# list of recorder names/emails is re-ordered into ordinal (case-sensitive,
# byte-value) sorted order to match the source system's canonical form.

function OrdinalLess($ordLeft, $ordRight) {
    $ordLenLeft = $ordLeft.Length
    $ordLenRight = $ordRight.Length
    $ordMinLen = $ordLenLeft
    if ($ordLenRight -lt $ordMinLen) { $ordMinLen = $ordLenRight }

    $ordPos = 0
    $ordResult = $ordLenLeft -lt $ordLenRight
    while ($ordPos -lt $ordMinLen) {
        $ordCodeLeft = [int][char]$ordLeft.Substring($ordPos,1)
        $ordCodeRight = [int][char]$ordRight.Substring($ordPos,1)
        if ($ordCodeLeft -ne $ordCodeRight) {
            $ordResult = $ordCodeLeft -lt $ordCodeRight
            $ordPos = $ordMinLen
        } else {
            $ordPos = $ordPos + 1
        }
    }
    return $ordResult
}

function OrdinalSortParts($sortItems) {
    $sortArr = @($sortItems)
    $sortCount = $sortArr.Length
    $sortOuter = 1
    while ($sortOuter -lt $sortCount) {
        $sortKey = $sortArr[$sortOuter]
        $sortInner = $sortOuter - 1
        while ($sortInner -ge 0 -and (OrdinalLess $sortKey $sortArr[$sortInner])) {
            $sortArr[$sortInner+1] = $sortArr[$sortInner]
            $sortInner = $sortInner - 1
        }
        $sortArr[$sortInner+1] = $sortKey
        $sortOuter = $sortOuter + 1
    }
    return $sortArr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRowNum = $ws.UsedRange.Rows.Count

$rowNum = 2
while ($rowNum -le $lastRowNum) {
    $recordedByCell = $ws.Cells.Item($rowNum, 7)
    $recordedByVal = $recordedByCell.Value2

    if ($recordedByVal -ne $null -and $recordedByVal -ne "") {
        $nameParts = $recordedByVal -split ", "
        $sortedParts = OrdinalSortParts $nameParts
        $newRecordedByVal = $sortedParts -join ", "

        # NOTE: this runtime's `-eq`/`-ne` string comparisons are
        # case-insensitive even with `-c...` variants, so a guard like
        # `if ($newRecordedByVal -ne $recordedByVal)` would wrongly skip
        # cells whose only change is a letter-case swap (e.g. "system" vs
        # "System"). Always write the recomputed value back instead;
        # it is a correctness no-op when nothing actually changed.
        $recordedByCell.Value2 = $newRecordedByVal
    }

    $rowNum = $rowNum + 1
}
